# Add keyboard tracker and SQLite
# Append new cash-register rows to the "Наличные" sheet (rows 29-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Наличные")

$rows = @(
    @(7676096317, "печать (1).docx", 2, 0.4, "2025-06-26 07:59:21"),
    @(7676096317, "ee71121c452ebe95882145ee5a20077b.pdf", 1, 0.2, "2025-06-26 12:26:32"),
    @(7676096317, "b40d26a2f26fd1f152f53ad626929df7.pdf", 1, 0.2, "2025-06-26 13:12:00"),
    @(746382370, "Курсовая Прашкович.docx", 1, 0.2, "2025-06-26 14:59:04"),
    @(7676096317, "печать.pdf", 14, 2.8, "2025-06-26 19:56:01")
)

$startRow = 29
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
